$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Asl Sorveglianza": append the new week block (17/01/2022 - 23/01/2022)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Asl Sorveglianza")

$ws1.Cells.Item(101,1).Value = "17/01/2022 - 23/01/2022"
$ws1.Cells.Item(101,2).Value = "AZIENDA USL TOSCANA SUD-EST"
$ws1.Cells.Item(101,3).Value = 56

$ws1.Cells.Item(102,1).Value = "17/01/2022 - 23/01/2022"
$ws1.Cells.Item(102,2).Value = "AZIENDA USL TOSCANA CENTRO"
$ws1.Cells.Item(102,3).Value = 123

$ws1.Cells.Item(103,1).Value = "17/01/2022 - 23/01/2022"
$ws1.Cells.Item(103,2).Value = "AZIENDA USL TOSCANA NORD-OVEST"
$ws1.Cells.Item(103,2).Font.Color = 0
$ws1.Cells.Item(103,3).Value = 123

$ws1.Cells.Item(104,1).Value = "17/01/2022 - 23/01/2022"
$ws1.Cells.Item(104,2).Value = "ASL TA"
$ws1.Cells.Item(104,3).Value = 1

$ws1.Cells.Item(105,2).Value = "Totale"
$ws1.Cells.Item(105,3).Value = 303

# ---------------------------------------------------------------------
# Sheet "Professione": append the new week block (17/01/2022 - 23/01/2022)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Professione")

$ws2.Cells.Item(74,1).Value = "17/01/2022 - 23/01/2022"
$ws2.Cells.Item(74,2).Value = "Insegnante"
$ws2.Cells.Item(74,3).Value = 285
$ws2.Cells.Item(74,4).Value = 3297

$ws2.Cells.Item(75,1).Value = "17/01/2022 - 23/01/2022"
$ws2.Cells.Item(75,2).Value = "Personale non docente"
$ws2.Cells.Item(75,3).Value = 18
$ws2.Cells.Item(75,4).Value = 162

$ws2.Cells.Item(76,2).Value = "Totale"
$ws2.Cells.Item(76,3).Value = 303
$ws2.Cells.Item(76,4).Value = 3459

# Apply the thousands-separator number format to the whole "Totale contatti"
# column (this also stamps the blank separator rows between week blocks with
# an empty, styled D cell, matching the authored workbook).
$ws2.Range("D2:D76").NumberFormat = "#,##0"

# ---------------------------------------------------------------------
# Sheet "Sesso ed età": append the new week block (17/01/2022 - 23/01/2022)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sesso ed età")

$ws3.Cells.Item(129,1).Value = "17/01/2022 - 23/01/2022"
$ws3.Cells.Item(129,2).Value = "0-18"
$ws3.Cells.Item(129,3).Value = "F"
$ws3.Cells.Item(129,4).Value = 2

$ws3.Cells.Item(130,1).Value = "17/01/2022 - 23/01/2022"
$ws3.Cells.Item(130,2).Value = "19-34"
$ws3.Cells.Item(130,3).Value = "F"
$ws3.Cells.Item(130,4).Value = 46

$ws3.Cells.Item(131,1).Value = "17/01/2022 - 23/01/2022"
$ws3.Cells.Item(131,2).Value = "19-34"
$ws3.Cells.Item(131,2).Font.Color = 0
$ws3.Cells.Item(131,3).Value = "M"
$ws3.Cells.Item(131,4).Value = 5

$ws3.Cells.Item(132,1).Value = "17/01/2022 - 23/01/2022"
$ws3.Cells.Item(132,2).Value = "35-49"
$ws3.Cells.Item(132,3).Value = "F"
$ws3.Cells.Item(132,4).Value = 152

$ws3.Cells.Item(133,1).Value = "17/01/2022 - 23/01/2022"
$ws3.Cells.Item(133,2).Value = "35-49"
$ws3.Cells.Item(133,3).Value = "M"
$ws3.Cells.Item(133,4).Value = 8

$ws3.Cells.Item(134,1).Value = "17/01/2022 - 23/01/2022"
$ws3.Cells.Item(134,2).Value = "50-64"
$ws3.Cells.Item(134,3).Value = "F"
$ws3.Cells.Item(134,4).Value = 80

$ws3.Cells.Item(135,1).Value = "17/01/2022 - 23/01/2022"
$ws3.Cells.Item(135,2).Value = "50-64"
$ws3.Cells.Item(135,3).Value = "M"
$ws3.Cells.Item(135,4).Value = 5

$ws3.Cells.Item(136,1).Value = "17/01/2022 - 23/01/2022"
$ws3.Cells.Item(136,2).Value = "65-79"
$ws3.Cells.Item(136,3).Value = "F"
$ws3.Cells.Item(136,4).Value = 3

$ws3.Cells.Item(137,1).Value = "17/01/2022 - 23/01/2022"
$ws3.Cells.Item(137,2).Value = "65-79"
$ws3.Cells.Item(137,3).Value = "M"
$ws3.Cells.Item(137,4).Value = 2

# ---------------------------------------------------------------------
# View-state niceties (best effort; mirrors the scrolled/selected ranges
# recorded by Excel after the edit).
# ---------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("C106").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("D75").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("D129:D137").Select() | Out-Null
